# Auto-generated Excel COM-interop edit script
# Applies crypto price/volume refresh + row32/row33 swap per commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text cells (coin name / link columns B & C) ---
$textAddrs = @("B32", "C32", "B33", "C33")
$textVals  = @("EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near")
for ($i = 0; $i -lt $textAddrs.Length; $i++) {
    $ws.Range($textAddrs[$i]).Value = $textVals[$i]
}

# --- Price cells (column D) - force as text so values like "1.00" / "0.528" ---
# --- aren't silently reinterpreted as numbers by Excel's type inference.   ---
$priceAddrs = @("D2", "D3", "D5", "D6", "D7", "D9", "D13", "D15", "D16", "D17", "D18", "D20", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D31", "D32", "D33", "D35", "D38", "D39", "D42", "D43", "D44", "D45", "D46", "D48", "D50", "D51")
$priceVals  = @("69.963.50", "3.870.21", "611.25", "175.37", "3.867.14", "0.528", "40.02", "4.518.06", "3.869.98", "70.015.37", "7.48", "16.61", "506.65", "9.64", "0.742", "86.13", "0.0000143", "12.67", "10.56", "2.99", "32.98", "7.96", "1.00", "0.142", "478.38", "49.78", "2.97", "43.39", "8.55", "2.942.57", "140.13", "27.07", "2.43")
for ($i = 0; $i -lt $priceAddrs.Length; $i++) {
    $c = $ws.Range($priceAddrs[$i])
    $c.NumberFormat = "@"
    $c.Value = $priceVals[$i]
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

# --- Volume(1h) percentage cells (column E) ---
$volAddrs = @("E2", "E3", "E4", "E5", "E6", "E7", "E8", "E9", "E10", "E11", "E12", "E13", "E14", "E15", "E16", "E17", "E18", "E19", "E20", "E21", "E22", "E23", "E24", "E25", "E26", "E27", "E28", "E29", "E30", "E31", "E32", "E33", "E34", "E35", "E36", "E37", "E38", "E39", "E40", "E41", "E42", "E43", "E44", "E45", "E46", "E47", "E48", "E50", "E51")
$volVals  = @("  -0.10%  ", "  +3.70%  ", "  +0.00%  ", "  -1.57%  ", "  -2.83%  ", "  +3.66%  ", "  +0.03%  ", "  -1.08%  ", "  +0.03%  ", "  +2.76%  ", "  -0.75%  ", "  -1.71%  ", "  -1.50%  ", "  +3.72%  ", "  +3.64%  ", "  -0.08%  ", "  -1.70%  ", "  -3.04%  ", "  -0.72%  ", "  +0.25%  ", "  +3.83%  ", "  +2.85%  ", "  -3.12%  ", "  -0.61%  ", "  +4.90%  ", "  -3.26%  ", "  -8.04%  ", "  +0.24%  ", "  +2.62%  ", "  +2.39%  ", "  +5.80%  ", "  +0.44%  ", "  -1.33%  ", "  +0.00%  ", "  -1.14%  ", "  -0.51%  ", "  +2.53%  ", "  +10.15%  ", "  -0.54%  ", "  -0.87%  ", "  -0.80%  ", "  +2.93%  ", "  -4.61%  ", "  -1.61%  ", "  -1.94%  ", "  -0.22%  ", "  +2.17%  ", "  -1.48%  ", "  -2.66%  ")
for ($i = 0; $i -lt $volAddrs.Length; $i++) {
    $ws.Range($volAddrs[$i]).Value = $volVals[$i]
}
